# fall 23 week 1 inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 1.4
$ws.Range("C4").Value = 1.45
$ws.Range("E4").Value = 1.25
$ws.Range("D5").Value = 1.32
$ws.Range("F5").Value = 1.03
$ws.Range("G5").Value = 0.67
$ws.Range("E6").Value = 1.34
$ws.Range("G6").Value = 1.03
$ws.Range("E7").Value = 1.97
